$d = $word.ActiveDocument

function Get-ParagraphByText($doc, [string]$needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text -like "*$needle*") {
            return $p
        }
    }
    return $null
}

# --- Paragraph 1: "... protected by a lock" pseudocode comment ---
$p1 = Get-ParagraphByText $d "Finding whether the accessed variable is protected by a lock"
if ($null -eq $p1) { throw "paragraph 1 not found" }
$xml1 = '<w:p w14:paraId="31051B8E" w14:textId="15A6C8C8" w:rsidR="00786BD3" w:rsidRPr="00786BD3" w:rsidRDefault="00786BD3" w:rsidP="004E20C5"><w:pPr><w:spacing w:line="240" w:lineRule="exact"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:i/><w:iCs/><w:color w:val="0070C0"/></w:rPr></w:pPr><w:r w:rsidRPr="00786BD3"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" w:hint="eastAsia"/><w:i/><w:iCs/><w:color w:val="0070C0"/></w:rPr><w:t>/</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:i/><w:iCs/><w:color w:val="0070C0"/></w:rPr><w:t xml:space="preserve">/ </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:i/><w:iCs/><w:color w:val="0070C0"/></w:rPr><w:t>Checking</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:i/><w:iCs/><w:color w:val="0070C0"/></w:rPr><w:t xml:space="preserve"> whether the accessed </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:i/><w:iCs/><w:color w:val="0070C0"/></w:rPr><w:t>f</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" w:hint="eastAsia"/><w:i/><w:iCs/><w:color w:val="0070C0"/></w:rPr><w:t>ield</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:i/><w:iCs/><w:color w:val="0070C0"/></w:rPr><w:t xml:space="preserve"> is protected by </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:i/><w:iCs/><w:color w:val="0070C0"/></w:rPr><w:t>any</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:i/><w:iCs/><w:color w:val="0070C0"/></w:rPr><w:t xml:space="preserve"> lock</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:i/><w:iCs/><w:color w:val="0070C0"/></w:rPr><w:t>.</w:t></w:r></w:p>'
$p1.Range.InsertXML($xml1)

# --- Paragraph 2: "... var exist in a key field" pseudocode comment ---
$p2 = Get-ParagraphByText $d "Finding whether the accessed var exist in a key field"
if ($null -eq $p2) { throw "paragraph 2 not found" }
$xml2 = '<w:p w14:paraId="450853E8" w14:textId="1CC6E75E" w:rsidR="00786BD3" w:rsidRPr="00786BD3" w:rsidRDefault="00786BD3" w:rsidP="00D708E1"><w:pPr><w:spacing w:line="240" w:lineRule="exact"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:i/><w:iCs/><w:color w:val="0070C0"/></w:rPr></w:pPr><w:r w:rsidRPr="00786BD3"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" w:hint="eastAsia"/><w:i/><w:iCs/><w:color w:val="0070C0"/></w:rPr><w:t>/</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:i/><w:iCs/><w:color w:val="0070C0"/></w:rPr><w:t xml:space="preserve">/ </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:i/><w:iCs/><w:color w:val="0070C0"/></w:rPr><w:t>Checking</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:i/><w:iCs/><w:color w:val="0070C0"/></w:rPr><w:t xml:space="preserve"> whether the accessed var exist in a</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:i/><w:iCs/><w:color w:val="0070C0"/></w:rPr><w:t>ny</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:i/><w:iCs/><w:color w:val="0070C0"/></w:rPr><w:t xml:space="preserve"> key field</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:i/><w:iCs/><w:color w:val="0070C0"/></w:rPr><w:t>.</w:t></w:r></w:p>'
$p2.Range.InsertXML($xml2)

Write-Host "Done."
